$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The original row 13 (B13/C13 held "5840535 - Messias Borges Silva", with no A13 label)
# is removed entirely; this shifts every subsequent row up by one, row 24 -> row 23.
$ws.Rows.Item(13).Delete()

# Row 10 (Objetivos:) now shows the docente's name instead of the long objectives text.
$ws.Range("B10:C10").Value = "5840535 - Messias Borges Silva"

# Row 13 (now "Programa resumido:") gets the short "Semestral" value.
$ws.Range("B13:C13").Value = "Semestral"

# Row 15 (now "Programa:") re-uses the activation-date text verbatim (copy, so it stays
# text rather than being reinterpreted as a date serial number).
$ws.Range("B8").Copy($ws.Range("B15"))
$ws.Range("C8").Copy($ws.Range("C15"))

# Row 18 (now "Método:") shows the docente's name again.
$ws.Range("B18:C18").Value = "5840535 - Messias Borges Silva"

# Row 19 (now "Critério:") gets "duas provas escritas".
$ws.Range("B19:C19").Value = "duas provas escritas"

# Row 20 (now "Norma de recuperação:") gets the evaluation-criteria paragraph.
$ws.Range("B20:C20").Value = "serão avaliados os conteúdos discutidos em sala e constantes da ementa do curso.A média da disciplina será a média aritmética das duas provas."

# Row 21 (now "Bibliografia:") gets the recovery-exam text.
$ws.Range("B21:C21").Value = "uma prova escrita com conteúdo de todo o semestre"
